$d = $word.ActiveDocument

# 1. Fix double .htm.htm extension for odds-calculator
$d.Content.Find.Execute(
    "odds-calculator.htm.htm", $true, $false, $false, $false, $false,
    $true, 1, $false, "odds-calculator.htm", 2)

# 2. Fix double .htm.htm extension for parlay-calculator
$d.Content.Find.Execute(
    "parlay-calculator.htm.htm", $true, $false, $false, $false, $false,
    $true, 1, $false, "parlay-calculator.htm", 2)

# 3. Replace competitor brand link: bet365-review -> 22bet-review
$d.Content.Find.Execute(
    "/sport/betting/ireland/bet365-review.htm", $true, $false, $false, $false, $false,
    $true, 1, $false, "/sport/betting/ireland/22bet-review.htm", 2)

# 4. Replace competitor brand link: paddy-power-review -> lunubet-review
$d.Content.Find.Execute(
    "/sport/betting/ireland/paddy-power-review.htm", $true, $false, $false, $false, $false,
    $true, 1, $false, "/sport/betting/ireland/lunubet-review.htm", 2)
